# Re-create the two lost stations ("Jan Pieter Heijestraat" and "Ten
# Katestraat") that the previous retrieve-function edit had dropped,
# re-inserting them between "Witte de Withstraat" and "Bilderdijkstraat".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Nicolaas Beetsstraat" (old row 13) is replaced by two rows: insert a
# fresh row below it so both new station names have a home, then
# overwrite the values in place.
$ws.Rows.Item(14).Insert()
$ws.Range("A13").Value = "Jan Pieter Heijestraat"
$ws.Range("A14").Value = "Ten Katestraat"

# Reflect the editor's last selection/cursor position.
$ws.Range("C15").Select()
